$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 489-490; this pushes the existing rows 489-498 down
# to 491-500 (Excel copies the row-above formatting, including the date
# style on column D, automatically).
$ws.Rows("489:490").Insert()

# --- Row 489 : new weekly "Primera" price point -----------------------
$ws.Range("A489").Value = 6
$ws.Range("B489").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C489").Value = "Metropolitana"
$ws.Range("D489").Value = 44448
$ws.Range("E489").Value = 13
$ws.Range("F489").Value = 100112040
$ws.Range("G489").Value = "Cilantro"
$ws.Range("H489").Value = "Sin especificar"
$ws.Range("I489").Value = "Primera"
$ws.Range("J489").Value = 650
$ws.Range("K489").Value = 4500
$ws.Range("L489").Value = 5000
$ws.Range("M489").Value = 4723
$ws.Range("N489").Value = "$/caja 36 atados"
$ws.Range("O489").Value = "Región Metropolitana"
$ws.Range("P489").Value = 131
$ws.Range("Q489").Value = 36
$ws.Range("R489").Value = "Hortaliza"

# --- Row 490 : new weekly "Primera" price point -----------------------
$ws.Range("A490").Value = 6
$ws.Range("B490").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C490").Value = "Metropolitana"
$ws.Range("D490").Value = 44448
$ws.Range("E490").Value = 13
$ws.Range("F490").Value = 100112040
$ws.Range("G490").Value = "Cilantro"
$ws.Range("H490").Value = "Sin especificar"
$ws.Range("I490").Value = "Primera"
$ws.Range("J490").Value = 510
$ws.Range("K490").Value = 8000
$ws.Range("L490").Value = 9000
$ws.Range("M490").Value = 8431
$ws.Range("N490").Value = "$/docena de atados"
$ws.Range("O490").Value = "Región Metropolitana"
$ws.Range("P490").Value = 2810
$ws.Range("Q490").Value = 3
$ws.Range("R490").Value = "Hortaliza"
